$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.682.31"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.524.75"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.72%  "
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "2.912.77"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "2.545.00"
$ws.Range("E16").Value = "  -3.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("D18").Value = "42.658.60"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.75%  "
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  +10.58%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("E37").Value = "  -7.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.27%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "1.990.29"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "2.767.31"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.15%  "
